$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

$ws.Range("E3").Value = 10.7
$ws.Range("E4").Value = 10.67
$ws.Range("F4").Value = 9.98
$ws.Range("C5").Value = 9.26
$ws.Range("D5").Value = 9.33
$ws.Range("H5").Value = 8.609999999999999
$ws.Range("D6").Value = 10.02
$ws.Range("I6").Value = 8.710000000000001
$ws.Range("H7").Value = 10.1
$ws.Range("J7").Value = 9.789999999999999
$ws.Range("E8").Value = 11.39
$ws.Range("G8").Value = 9.9
$ws.Range("F9").Value = 11.29
$ws.Range("G10").Value = 10.21
